$wb = $excel.ActiveWorkbook

# Add a new worksheet at the end of the workbook and name it "Relation"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Relation"

# Column widths (engine's char-width unit differs slightly from Excel's MDW-based
# unit, so these are the closest achievable approximation of the authored widths
# 31.59765625 and 69.265625 for columns B and C respectively)
$newSheet.Columns.Item(2).ColumnWidth = 30.833333333333332
$newSheet.Columns.Item(3).ColumnWidth = 68.5

$newSheet.Range("B4").Value = 'Rules Integer indexInt_Relation_ParamDouble(Double param)'
$newSheet.Range("B5").Value = 'C1'
$newSheet.Range("C5").Value = 'RET1'
$newSheet.Range("B6").Value = 'param < x'
$newSheet.Range("B7").Value = 'Integer x'
$newSheet.Range("B8").Value = 'condition'
$newSheet.Range("B9").Value = 12
$newSheet.Range("C9").Value = 1
$newSheet.Range("C10").Value = 0
$newSheet.Range("C11").Value = 0
$newSheet.Range("B13").Value = 'Test indexInt_Relation_ParamDouble indexInt_Relation_ParamDoubleTest'
$newSheet.Range("B14").Value = 'param'
$newSheet.Range("C14").Value = '_res_'
$newSheet.Range("B15").Value = 'param'
$newSheet.Range("C15").Value = 'Result'
$newSheet.Range("B16").Value = 11.99
$newSheet.Range("C16").Value = 1
$newSheet.Range("B17").Value = 12.1
$newSheet.Range("C17").Value = 0
$newSheet.Range("B21").Value = 'Rules Integer indexDouble_Relation_ParamInteger(Integer param)'
$newSheet.Range("B22").Value = 'C1'
$newSheet.Range("C22").Value = 'RET1'
$newSheet.Range("B23").Value = 'param < x'
$newSheet.Range("B24").Value = 'Double x'
$newSheet.Range("B25").Value = 'condition'
$newSheet.Range("B26").Value = 12.1
$newSheet.Range("C26").Value = 1
$newSheet.Range("C27").Value = 0
$newSheet.Range("C28").Value = 0
$newSheet.Range("B31").Value = 'Test  indexDouble_Relation_ParamInteger indexDouble_Relation_ParamIntegerTest'
$newSheet.Range("B32").Value = 'param'
$newSheet.Range("C32").Value = '_res_'
$newSheet.Range("B33").Value = 'Param'
$newSheet.Range("C33").Value = '_res_'
$newSheet.Range("B34").Value = 12
$newSheet.Range("C34").Value = 1
$newSheet.Range("B35").Value = 13
$newSheet.Range("C35").Value = 0

# empty cell with quote-prefix style, matching original artifact at B59
$newSheet.Range("B59").Value = "'"
$newSheet.Range("B59").Value = ""

# Restore selection to match target view state
$newSheet.Range("C27").Select() | Out-Null
